# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header already labeled "K") is recalculated
# and rewritten with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 4
    4  = 2
    5  = 1
    6  = 3
    7  = 1
    8  = 1
    9  = 1
    10 = 6
    11 = 1
    12 = 3
    13 = 4
    14 = 1
    15 = 2
    16 = 5
    17 = 0
    18 = 7
    19 = 4
    20 = 5
    21 = 2
    22 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
